$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "presentation figures and finer time steps" -- refresh the dilution/value
# table (B2:O10) with the new measurements. Row 1 is the header row and is
# untouched; only individual measurement cells in rows 2-10 change.

# Row 2
$ws.Range("B2").Value = 25
$ws.Range("D2").Value = 24
$ws.Range("F2").Value = 12
$ws.Range("H2").Value = 10
$ws.Range("J2").Value = 5
$ws.Range("L2").Value = 29
$ws.Range("N2").Value = 28
$ws.Range("O2").Value = 6

# Row 3
$ws.Range("B3").Value = 26
$ws.Range("D3").Value = 15
$ws.Range("F3").Value = 8
$ws.Range("H3").Value = 6
$ws.Range("I3").Value = 5
$ws.Range("J3").Value = 25
$ws.Range("L3").Value = 18
$ws.Range("N3").Value = 21
$ws.Range("O3").Value = 6

# Row 4
$ws.Range("B4").Value = 19
$ws.Range("D4").Value = 15
$ws.Range("F4").Value = 39
$ws.Range("H4").Value = 28
$ws.Range("J4").Value = 15
$ws.Range("L4").Value = 22
$ws.Range("N4").Value = 31
$ws.Range("O4").Value = 6

# Row 5
$ws.Range("B5").Value = 24
$ws.Range("D5").Value = 11
$ws.Range("F5").Value = 44
$ws.Range("H5").Value = 28
$ws.Range("J5").Value = 14
$ws.Range("L5").Value = 11
$ws.Range("N5").Value = 18

# Row 6
$ws.Range("B6").Value = 25
$ws.Range("F6").Value = 31
$ws.Range("H6").Value = 11
$ws.Range("J6").Value = 36
$ws.Range("L6").Value = 6
$ws.Range("M6").Value = 6
$ws.Range("N6").Value = 17

# Row 7
$ws.Range("B7").Value = 25
$ws.Range("D7").Value = 7
$ws.Range("F7").Value = 23
$ws.Range("H7").Value = 34
$ws.Range("J7").Value = 10
$ws.Range("L7").Value = 38
$ws.Range("N7").Value = 18

# Row 8
$ws.Range("B8").Value = 25
$ws.Range("D8").Value = 8
$ws.Range("F8").Value = 12
$ws.Range("H8").Value = 25
$ws.Range("J8").Value = 39
$ws.Range("L8").Value = 10
$ws.Range("M8").Value = 5
$ws.Range("N8").Value = 17

# Row 9
$ws.Range("D9").Value = 6
$ws.Range("F9").Value = 8
$ws.Range("H9").Value = 10
$ws.Range("J9").Value = 24
$ws.Range("L9").Value = 16
$ws.Range("M9").Value = 4
$ws.Range("N9").Value = 16

# Row 10
$ws.Range("B10").Value = 19
$ws.Range("D10").Value = 7
$ws.Range("F10").Value = 8
$ws.Range("G10").Value = 3
$ws.Range("H10").Value = 13
$ws.Range("J10").Value = 9
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 20
$ws.Range("N10").Value = 13

# View state left behind by the author: scrolled right one column, with
# M9 as the active selection (previously the freshly-appended O11).
$win = $excel.ActiveWindow
$win.ScrollColumn = 2
$win.ScrollRow = 1
$ws.Range("M9").Select()
